$d = $word.ActiveDocument

# Helper: find the InlineShape that lives inside a header/footer Range and
# rename it. Going through HeaderFooter.Range.InlineShapes(n) directly can
# address a stale/ambiguous block when the picture isn't in the first
# paragraph of the story, so we walk the paragraphs and grab the shape from
# the paragraph Range that actually contains it - that handle resolves
# reliably for both Set and Get.
function Rename-StoryPictures($story, [string]$newName) {
    if (-not $story.Exists) {
        return
    }
    $paras = $story.Range.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $paraRange = $paras.Item($i).Range
        $shapes = $paraRange.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shapes.Item($j).Name = $newName
        }
    }
}

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections($s)

    for ($h = 1; $h -le $section.Headers.Count; $h++) {
        Rename-StoryPictures $section.Headers($h) "image2.jpg"
    }

    for ($f = 1; $f -le $section.Footers.Count; $f++) {
        Rename-StoryPictures $section.Footers($f) "image1.png"
    }
}
